# Update the "Expenses" monthly report sheet (TDSheet) to:
#  - correct the existing "Апрель 2023 г." row (row 137) totals
#  - append a new row for "Май 2023 г." (row 138)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TDSheet")

# Row 137 ("Апрель 2023 г.") totals changed.
$ws.Cells.Item(137, 2).Value = 149767.79999999999
$ws.Cells.Item(137, 3).Value = 14740.002

# New row 138 for "Май 2023 г.": copy the formatting (styles + row height) of
# the row above it down into the freshly added row, then fill in its values.
$ws.Range("A137:C137").Copy()
$ws.Range("A138:C138").PasteSpecial(-4122)
$ws.Rows.Item(138).RowHeight = $ws.Rows.Item(137).RowHeight

$ws.Cells.Item(138, 1).Value = "Май 2023 г."
$ws.Cells.Item(138, 2).Value = 63297.4
$ws.Cells.Item(138, 3).Value = 7371.4

# Move the active selection to the next empty cell in column C, as Excel would
# after the user finished entering the new row.
$ws.Range("C139").Select()
